$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-15 (rows 16-19 stay unchanged)
$data = @(
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Al Horford", "PF,C", "Boston Celtics"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Jakob Poeltl", "C", "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
